$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1755.5769
$ws.Range("I40").Value = 1277.8125
$ws.Range("J40").Value = 2520
$ws.Range("K40").Value = 1277.8125
$ws.Range("L40").Value = 2520
$ws.Range("M40").Value = -1102.8125
$ws.Range("N40").Value = -2870
$ws.Range("H138").Value = 7865
$ws.Range("I138").Value = 2266.6667
$ws.Range("J138").Value = 9964.375
$ws.Range("K138").Value = 6800.000100000001
$ws.Range("L138").Value = 29893.125
$ws.Range("M138").Value = -1660.000100000001
$ws.Range("N138").Value = -40173.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8835
$ws.Range("I2").Value = 30011
$ws.Range("J2").Value = 4599.8
$ws.Range("K2").Value = 30011
$ws.Range("L2").Value = 4599.8
$ws.Range("M2").Value = -29898
$ws.Range("N2").Value = -4825.8
$ws.Range("H45").Value = 1281.8572
$ws.Range("I45").Value = 1143.75
$ws.Range("J45").Value = 1466
$ws.Range("K45").Value = 1143.75
$ws.Range("L45").Value = 1466
$ws.Range("M45").Value = -766.75
$ws.Range("N45").Value = -2220
$ws.Range("H116").Value = 8835
$ws.Range("I116").Value = 30011
$ws.Range("J116").Value = 4599.8
$ws.Range("K116").Value = 30011
$ws.Range("L116").Value = 4599.8
$ws.Range("M116").Value = -27717
$ws.Range("N116").Value = -9187.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8835
$ws.Range("I3").Value = 30011
$ws.Range("J3").Value = 4599.8
$ws.Range("K3").Value = 30011
$ws.Range("L3").Value = 4599.8
$ws.Range("M3").Value = -29897
$ws.Range("N3").Value = -4827.8
$ws.Range("H105").Value = 7577277.5
$ws.Range("I105").Value = 7577277.5
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 7577277.5
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -7575530.5
$ws.Range("H134").Value = 1807.8966
$ws.Range("I134").Value = 1514.5264
$ws.Range("J134").Value = 2365.3
$ws.Range("K134").Value = 4543.5792
$ws.Range("L134").Value = 7095.900000000001
$ws.Range("M134").Value = -2008.5792
$ws.Range("N134").Value = -12165.9
$ws.Range("N105").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3313.2727
$ws.Range("I16").Value = 2006.5714
$ws.Range("J16").Value = 5600
$ws.Range("K16").Value = 2006.5714
$ws.Range("L16").Value = 5600
$ws.Range("M16").Value = -1719.5714
$ws.Range("N16").Value = -6174
$ws.Range("H41").Value = 13900
$ws.Range("I41").Value = 8000
$ws.Range("J41").Value = 19800
$ws.Range("K41").Value = 8000
$ws.Range("L41").Value = 19800
$ws.Range("M41").Value = -7572
$ws.Range("N41").Value = -20656
$ws.Range("H50").Value = 14350
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 14350
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 14350
$ws.Range("N50").Value = -15600
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("H68").Value = 47995
$ws.Range("J68").Value = 47995
$ws.Range("L68").Value = 47995
$ws.Range("N68").Value = -49493
$ws.Range("H71").Value = 47995
$ws.Range("J71").Value = 47995
$ws.Range("L71").Value = 143985
$ws.Range("N71").Value = -151473
$ws.Range("H74").Value = 250022160
$ws.Range("J74").Value = 250022160
$ws.Range("L74").Value = 250022160
$ws.Range("N74").Value = -250023908
$ws.Range("H77").Value = 250022160
$ws.Range("J77").Value = 250022160
$ws.Range("L77").Value = 750066480
$ws.Range("N77").Value = -750075216
$ws.Range("H113").Value = 3313.2727
$ws.Range("I113").Value = 2006.5714
$ws.Range("J113").Value = 5600
$ws.Range("K113").Value = 2006.5714
$ws.Range("L113").Value = 5600
$ws.Range("M113").Value = 163.4286
$ws.Range("N113").Value = -9940
$ws.Range("M50").Value = ""
$ws.Range("N51").Value = ""
$ws.Range("M51").Value = ""
$ws.Range("N61").Value = ""
$ws.Range("M61").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 8548.416999999999
$ws.Range("I2").Value = 48.75
$ws.Range("J2").Value = 12798.25
$ws.Range("K2").Value = 292.5
$ws.Range("L2").Value = 76789.5
$ws.Range("M2").Value = -179.5
$ws.Range("N2").Value = -77015.5
$ws.Range("H38").Value = 117.70588
$ws.Range("I38").Value = 80.90909000000001
$ws.Range("J38").Value = 185.16667
$ws.Range("K38").Value = 242.72727
$ws.Range("L38").Value = 555.50001
$ws.Range("M38").Value = 104.27273
$ws.Range("N38").Value = -1249.50001
$ws.Range("H68").Value = 1379.4884
$ws.Range("I68").Value = 775.75
$ws.Range("J68").Value = 2142.1052
$ws.Range("K68").Value = 2327.25
$ws.Range("L68").Value = 6426.3156
$ws.Range("M68").Value = -1516.25
$ws.Range("N68").Value = -8048.3156
$ws.Range("H71").Value = 1379.4884
$ws.Range("I71").Value = 775.75
$ws.Range("J71").Value = 2142.1052
$ws.Range("K71").Value = 6981.75
$ws.Range("L71").Value = 19278.9468
$ws.Range("M71").Value = -2925.75
$ws.Range("N71").Value = -27390.9468
$ws.Range("H107").Value = 349.20456
$ws.Range("I107").Value = 226.97223
$ws.Range("J107").Value = 899.25
$ws.Range("K107").Value = 680.91669
$ws.Range("L107").Value = 2697.75
$ws.Range("M107").Value = 1239.08331
$ws.Range("N107").Value = -6537.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("H132").Value = 8280.666999999999
$ws.Range("I132").Value = 5576.75
$ws.Range("J132").Value = 23422.6
$ws.Range("K132").Value = 16730.25
$ws.Range("L132").Value = 70267.79999999999
$ws.Range("M132").Value = -14200.25
$ws.Range("N132").Value = -75327.79999999999
$ws.Range("M58").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2437.8147
$ws.Range("I7").Value = 2076.4285
$ws.Range("J7").Value = 3702.6667
$ws.Range("K7").Value = 2076.4285
$ws.Range("L7").Value = 3702.6667
$ws.Range("M7").Value = -1964.4285
$ws.Range("N7").Value = -3926.6667
$ws.Range("H16").Value = 2489
$ws.Range("I16").Value = 2489
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2489
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2319
$ws.Range("H68").Value = 10025.167
$ws.Range("I68").Value = 13787.75
$ws.Range("K68").Value = 13787.75
$ws.Range("M68").Value = -13038.75
$ws.Range("H71").Value = 10025.167
$ws.Range("I71").Value = 13787.75
$ws.Range("K71").Value = 68938.75
$ws.Range("M71").Value = -65194.75
$ws.Range("H126").Value = 2437.8147
$ws.Range("I126").Value = 2076.4285
$ws.Range("J126").Value = 3702.6667
$ws.Range("K126").Value = 6229.2855
$ws.Range("L126").Value = 11108.0001
$ws.Range("M126").Value = -3759.2855
$ws.Range("N126").Value = -16048.0001
$ws.Range("N16").Value = ""
